# Apply the "Allow policy mandated capacity additions to overwrite new power
# plant ban, and set ban to apply to all plants from 2021-2023" edit.

$wb = $excel.ActiveWorkbook

$wsBBNPPTY = $wb.Worksheets.Item("BBNPPTY")

# Make the BBNPPTY sheet the active (selected) tab instead of "About".
$wsBBNPPTY.Activate()

# The ban (Boolean 0/1 flag) now applies to every plant type for the years
# 2021-2023 (columns B, C, D), since historical additions for those years
# are already specified elsewhere in the input data.
$wsBBNPPTY.Range("B2:D25").Value = 1

# Update the visible selection to match where the edit was made.
$wsBBNPPTY.Range("B2:D25").Select()
